# LhARABeamLine-Params-LsrDrvn-Gabor.xlsx
#
# Commit: "Issues at low energy, worked around in particle"
#
# Two visible, scriptable changes:
#  1. The sheet tab was renamed, picking up the workbook's "-Ga" (Gabor)
#     suffix: "LhARABeamLine-Params-LsrDrvn" -> "LhARABeamLine-Params-LsrDrvn-Ga"
#  2. The view was scrolled back up (away from row 71, where the low-energy
#     parameters near the bottom of the sheet had been worked on) and the
#     cursor left resting on F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "LhARABeamLine-Params-LsrDrvn-Ga"

$ws.Range("F8").Select()
